$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "I, _KARTHIK SUBRAMANIAN RAMNATH" -> "I, _______________________________"
#    + "_____"  (blank the applicant's name, drop the underline formatting)
# ---------------------------------------------------------------------------

# 1a. Change the "I, _" run's text first (while the name run still carries its
#     own underline formatting) so it cannot be merged with the run that
#     follows it.
$find1a = $d.Content.Find
$find1a.ClearFormatting()
$find1a.Text = "I, _"
$find1a.Replacement.ClearFormatting()
$find1a.Replacement.Text = "I, _______________________________"
$find1a.Execute([ref]"I, _", $false, $false, $false, $false, $false, $true, 1, $false, `
    "I, _______________________________", 2) | Out-Null

# 1b. Now blank out the name itself and strip its underline formatting.
$find1b = $d.Content.Find
$find1b.ClearFormatting()
$find1b.Text = "KARTHIK SUBRAMANIAN RAMNATH"
$find1b.Replacement.ClearFormatting()
$find1b.Replacement.Font.Underline = 0
$find1b.Replacement.Text = "_____"
$find1b.Execute([ref]"KARTHIK SUBRAMANIAN RAMNATH", $false, $false, $false, $false, $false, $true, 1, $false, `
    "_____", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Identity Card No. _P6426635," -> "Identity Card No. ____________________________,"
#    merged into a single run with no underline.
# ---------------------------------------------------------------------------

# 2a. Replace the whole span (covering the three original runs) via a Range
#     assignment -- this merges them into one run (matching the formatting of
#     the run that follows) without touching the unrelated run before it.
$rng2 = $d.Content.Duplicate
$rng2.Find.ClearFormatting()
$rng2.Find.Text = "Identity Card No. _P6426635,"
$rng2.Find.Execute() | Out-Null
$rng2.Font.Underline = 0
$rng2.Text = "Identity Card No. ____________________________,"

# 2b. Tidy up: a plain Range.Text assignment leaves an explicit
#     <w:u w:val="none"/>; run it back through Find/Replace (same text) with
#     Replacement formatting cleared so the underline element disappears
#     entirely instead of being serialized as "none".
$find2b = $d.Content.Find
$find2b.ClearFormatting()
$find2b.Text = "Identity Card No. ____________________________,"
$find2b.Replacement.ClearFormatting()
$find2b.Replacement.Font.Underline = 0
$find2b.Replacement.Text = "Identity Card No. ____________________________,"
$find2b.Execute([ref]"Identity Card No. ____________________________,", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Identity Card No. ____________________________,", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Accepted by " + "AVA :" + " " -> "Accepted by AVA : " as a single run,
#    dropping the gramStart/gramEnd proofErr markers between them.
# ---------------------------------------------------------------------------

# 3a. Push a placeholder through so the three runs actually get merged (an
#     assignment that is textually a no-op doesn't trigger the merge).
$rng3 = $d.Content.Duplicate
$rng3.Find.ClearFormatting()
$rng3.Find.Text = "Accepted by AVA : "
$rng3.Find.Execute() | Out-Null
$rng3.Text = "Accepted by AVA : @@PLACEHOLDER@@"

# 3b. Swap the placeholder back out for the real text -- by now it's one run
#     so this just rewrites its text in place.
$rng3b = $d.Content.Duplicate
$rng3b.Find.ClearFormatting()
$rng3b.Find.Text = "Accepted by AVA : @@PLACEHOLDER@@"
$rng3b.Find.Execute() | Out-Null
$rng3b.Text = "Accepted by AVA : "

Write-Host "edits applied"
